$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers for the two new columns (BI=61, BJ=62) ---
$ws.Cells.Item(1, 61).Value = "EditLAT"
$ws.Cells.Item(1, 62).Value = "EditLong"

# --- Edited Lat / Long data for rows 2-10 ---
$ws.Cells.Item(2, 61).Value = 28.596306999999999
$ws.Cells.Item(2, 62).Value = -97.89425

$ws.Cells.Item(3, 61).Value = 28.325261999999999
$ws.Cells.Item(3, 62).Value = -97.570344000000006

$ws.Cells.Item(4, 61).Value = 28.384394
$ws.Cells.Item(4, 62).Value = -97.848106999999999

$ws.Cells.Item(5, 61).Value = 28.460977
$ws.Cells.Item(5, 62).Value = -97.661371000000003

$ws.Cells.Item(6, 61).Value = 28.420452999999998
$ws.Cells.Item(6, 62).Value = -97.755132000000003

$ws.Cells.Item(7, 61).Value = 28.403562000000001
$ws.Cells.Item(7, 62).Value = -97.759055000000004

$ws.Cells.Item(8, 61).Value = 28.397148000000001
$ws.Cells.Item(8, 62).Value = -97.734679

$ws.Cells.Item(9, 61).Value = 28.388449999999999
$ws.Cells.Item(9, 62).Value = -97.743106999999995

$ws.Cells.Item(10, 61).Value = 28.219092
$ws.Cells.Item(10, 62).Value = -97.679786000000007

# --- Column widths (bestFit-style) for BD/BE (56/57) and the new BI/BJ (61/62) ---
# The underlying engine adds a fixed 5/6-character padding on top of whatever
# ColumnWidth is assigned, so we back that padding out to land on the target
# stored widths from the workbook.
$ws.Columns.Item(56).ColumnWidth = 11 - 5/6
$ws.Columns.Item(57).ColumnWidth = 11.6640625 - 5/6
$ws.Columns.Item(61).ColumnWidth = 10 - 5/6
$ws.Columns.Item(62).ColumnWidth = 10.6640625 - 5/6

# --- Scroll the view over to show the new columns, then reselect ---
$ws.Application.ActiveWindow.ScrollColumn = 42
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("BC19").Select() | Out-Null

# --- Force a pageSetup element (portrait) to be written, matching the diff ---
$ws.PageSetup.Orientation = 1
